# Pharma_Society_Report.xlsx update
# - WVOS membership count (B5): 87 -> 187
# - DSCO membership count (B6): "104" (text) -> 1104 (number)
# - OSNJ membership count (B7): "649" (text) -> 1649 (number)
# - ESHOS membership count (B8): "200" (text) -> 1200 (number)
# - Active selection moves from D11 to F16

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5").Value2 = 187
$ws.Range("B6").Value2 = 1104
$ws.Range("B7").Value2 = 1649
$ws.Range("B8").Value2 = 1200

$ws.Range("F16").Select()
